$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.742.77'
$ws.Range("E2").Value = '  +0.34%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.603.08'
$ws.Range("E3").Value = '  +0.36%  '

# Row 4
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.94'
$ws.Range("E5").Value = '  +0.22%  '

# Row 6
$ws.Range("E6").Value = '  -0.34%  '

# Row 7
$ws.Range("E7").Value = '  +0.17%  '

# Row 8
$ws.Range("E8").Value = '  +0.23%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.247'
$ws.Range("E9").Value = '  +0.32%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.73'

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("E11").Value = '  +1.19%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.828.48'
$ws.Range("E12").Value = '  +0.29%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.597.91'
$ws.Range("E13").Value = '  +0.95%  '

# Row 14
$ws.Range("E14").Value = '  +0.58%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.522'
$ws.Range("E15").Value = '  -0.27%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.11'
$ws.Range("E16").Value = '  -0.08%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₃0741'
$ws.Range("E17").Value = '  +0.41%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '210.34'

# Row 19
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.01'
$ws.Range("E19").Value = '  +0.20%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.19'
$ws.Range("E20").Value = '  +2.59%  '

# Row 21
$ws.Range("E21").Value = '  -0.20%  '

# Row 22
$ws.Range("E22").Value = '  -2.61%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.01'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.72'
$ws.Range("E24").Value = '  -1.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.01'
$ws.Range("E25").Value = '  +0.19%  '

# Row 27
$ws.Range("E27").Value = '  -0.95%  '

# Row 28
$ws.Range("E28").Value = '  +0.72%  '

# Row 29
$ws.Range("E29").Value = '  -0.66%  '

# Row 30
$ws.Range("E30").Value = '  +0.10%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.27'
$ws.Range("E31").Value = '  +1.11%  '

# Row 32
$ws.Range("E32").Value = '  +1.05%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.294.29'
$ws.Range("E33").Value = '  +0.84%  '

# Row 34
$ws.Range("E34").Value = '  +0.75%  '

# Row 35
$ws.Range("E35").Value = '  +0.69%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.600'
$ws.Range("E36").Value = '  -3.16%  '

# Row 37
$ws.Range("E37").Value = '  +11.09%  '

# Row 38
$ws.Range("E38").Value = '  -0.18%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.832'
$ws.Range("E39").Value = '  -0.35%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.40'

# Row 41
$ws.Range("E41").Value = '  +0.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.785'
$ws.Range("E42").Value = '  +0.09%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.85'
$ws.Range("E43").Value = '  -1.73%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.739.53'
$ws.Range("E44").Value = '  +0.22%  '

# Row 45
$ws.Range("E45").Value = '  -0.11%  '

# Row 46
$ws.Range("E46").Value = '  -1.62%  '

# Row 47
$ws.Range("E47").Value = '  -0.36%  '

# Row 48
$ws.Range("E48").Value = '  +1.65%  '

# Row 49
$ws.Range("B49").Value = 'USDD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.18%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.42'
$ws.Range("E50").Value = '  +0.15%  '

# Row 51
$ws.Range("E51").Value = '  +0.97%  '
